$wb = $excel.ActiveWorkbook

# Sheet 1: ALC
$ws = $wb.Worksheets.Item(1)
$ws.Range("H2").Value = 1238
$ws.Range("I2").Value = 1263.3334
$ws.Range("K2").Value = 1263.3334
$ws.Range("M2").Value = -1150.3334
$ws.Range("H8").Value = 358.42856
$ws.Range("I8").Value = 182.2
$ws.Range("K8").Value = 546.5999999999999
$ws.Range("M8").Value = -407.5999999999999
$ws.Range("H18").Value = 743.75
$ws.Range("I18").Value = 800
$ws.Range("J18").Value = 350
$ws.Range("K18").Value = 800
$ws.Range("L18").Value = 350
$ws.Range("M18").Value = -516
$ws.Range("N18").Value = -918
$ws.Range("H33").Value = 51295.5
$ws.Range("I33").Value = 63931
$ws.Range("K33").Value = 63931
$ws.Range("M33").Value = -63702
$ws.Range("H40").Value = 6999.0835
$ws.Range("I40").Value = 6285.7144
$ws.Range("K40").Value = 6285.7144
$ws.Range("M40").Value = -6110.7144
$ws.Range("H43").Value = 4565.702
$ws.Range("J43").Value = 4123.4546
$ws.Range("L43").Value = 4123.4546
$ws.Range("N43").Value = -4261.4546
$ws.Range("H55").Value = 122.5
$ws.Range("I55").Value = 195.625
$ws.Range("J55").Value = 98.125
$ws.Range("K55").Value = 195.625
$ws.Range("L55").Value = 98.125
$ws.Range("M55").Value = 18.375
$ws.Range("N55").Value = -526.125
$ws.Range("H80").Value = 3371.476
$ws.Range("I80").Value = 1676.7778
$ws.Range("K80").Value = 5030.3334
$ws.Range("M80").Value = -4032.3334
$ws.Range("H82").Value = 6091.3335
$ws.Range("I82").Value = 5808.727
$ws.Range("K82").Value = 17426.181
$ws.Range("M82").Value = -17020.181
$ws.Range("H83").Value = 3371.476
$ws.Range("I83").Value = 1676.7778
$ws.Range("K83").Value = 15091.0002
$ws.Range("M83").Value = -10099.0002
$ws.Range("H85").Value = 6091.3335
$ws.Range("I85").Value = 5808.727
$ws.Range("K85").Value = 17426.181
$ws.Range("M85").Value = -16022.181
$ws.Range("H86").Value = 1956.375
$ws.Range("I86").Value = 2142
$ws.Range("J86").Value = 1399.5
$ws.Range("K86").Value = 2142
$ws.Range("L86").Value = 1399.5
$ws.Range("M86").Value = -1019
$ws.Range("N86").Value = -3645.5
$ws.Range("H89").Value = 1956.375
$ws.Range("I89").Value = 2142
$ws.Range("J89").Value = 1399.5
$ws.Range("K89").Value = 10710
$ws.Range("L89").Value = 6997.5
$ws.Range("M89").Value = -5094
$ws.Range("N89").Value = -18229.5
$ws.Range("H115").Value = 1905.1818
$ws.Range("I115").Value = 1905.1818
$ws.Range("K115").Value = 5715.5454
$ws.Range("M115").Value = -4148.5454
$ws.Range("H116").Value = 53424.43
$ws.Range("I116").Value = 128830.875
$ws.Range("J116").Value = 7020.4614
$ws.Range("K116").Value = 128830.875
$ws.Range("L116").Value = 7020.4614
$ws.Range("M116").Value = -125388.875
$ws.Range("N116").Value = -13904.4614
$ws.Range("H125").Value = 4805.3335
$ws.Range("I125").Value = 3962.3333
$ws.Range("J125").Value = 5226.8335
$ws.Range("K125").Value = 35660.9997
$ws.Range("L125").Value = 47041.5015
$ws.Range("M125").Value = -33200.9997
$ws.Range("N125").Value = -51961.5015
$ws.Range("H132").Value = 31371.172
$ws.Range("I132").Value = 35013.742
$ws.Range("J132").Value = 3141.25
$ws.Range("K132").Value = 105041.226
$ws.Range("L132").Value = 9423.75
$ws.Range("M132").Value = -102511.226
$ws.Range("N132").Value = -14483.75
$ws.Range("H135").Value = 900.8
$ws.Range("I135").Value = 900.8
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 8107.2
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -5572.2
$ws.Range("N135").ClearContents()
$ws.Range("H137").Value = 2280.8845
$ws.Range("I137").Value = 2041.6364
$ws.Range("J137").Value = 2456.3333
$ws.Range("K137").Value = 6124.9092
$ws.Range("L137").Value = 7368.999899999999
$ws.Range("M137").Value = -3574.9092
$ws.Range("N137").Value = -12468.9999
$ws.Range("H138").Value = 7031.263
$ws.Range("J138").Value = 7023.49
$ws.Range("L138").Value = 21070.47
$ws.Range("N138").Value = -31350.47
$ws.Range("H141").Value = 3554.0833
$ws.Range("I141").Value = 3554.0833
$ws.Range("K141").Value = 10662.2499
$ws.Range("M141").Value = -5482.249899999999

# Sheet 2: ARM
$ws = $wb.Worksheets.Item(2)
$ws.Range("H32").Value = 2617.239
$ws.Range("I32").Value = 2451.0698
$ws.Range("J32").Value = 4999
$ws.Range("K32").Value = 2451.0698
$ws.Range("L32").Value = 4999
$ws.Range("M32").Value = -2164.0698
$ws.Range("N32").Value = -5573
$ws.Range("H45").Value = 1351.7059
$ws.Range("I45").Value = 1365.3334
$ws.Range("J45").Value = 1249.5
$ws.Range("K45").Value = 1365.3334
$ws.Range("L45").Value = 1249.5
$ws.Range("M45").Value = -988.3334
$ws.Range("N45").Value = -2003.5
$ws.Range("H61").Value = 5999
$ws.Range("I61").Value = 4998.5
$ws.Range("J61").Value = 8000
$ws.Range("K61").Value = 4998.5
$ws.Range("L61").Value = 8000
$ws.Range("M61").Value = -4786.5
$ws.Range("N61").Value = -8424
$ws.Range("H63").Value = 1849
$ws.Range("I63").Value = 2130.4
$ws.Range("J63").Value = 1692.6666
$ws.Range("K63").Value = 2130.4
$ws.Range("L63").Value = 1692.6666
$ws.Range("M63").Value = -1444.4
$ws.Range("N63").Value = -3064.6666
$ws.Range("H66").Value = 1849
$ws.Range("I66").Value = 2130.4
$ws.Range("J66").Value = 1692.6666
$ws.Range("K66").Value = 10652
$ws.Range("L66").Value = 8463.333000000001
$ws.Range("M66").Value = -7220
$ws.Range("N66").Value = -15327.333
$ws.Range("H74").Value = 2958.4866
$ws.Range("I74").Value = 3013.3713
$ws.Range("J74").Value = 1998
$ws.Range("K74").Value = 3013.3713
$ws.Range("L74").Value = 1998
$ws.Range("M74").Value = -2139.3713
$ws.Range("N74").Value = -3746
$ws.Range("H77").Value = 2958.4866
$ws.Range("I77").Value = 3013.3713
$ws.Range("J77").Value = 1998
$ws.Range("K77").Value = 15066.8565
$ws.Range("L77").Value = 9990
$ws.Range("M77").Value = -10698.8565
$ws.Range("N77").Value = -18726
$ws.Range("H88").Value = 1964.8889
$ws.Range("I88").Value = 1848
$ws.Range("J88").Value = 2023.3334
$ws.Range("K88").Value = 1848
$ws.Range("L88").Value = 2023.3334
$ws.Range("M88").Value = -1442
$ws.Range("N88").Value = -2835.3334
$ws.Range("H91").Value = 1964.8889
$ws.Range("I91").Value = 1848
$ws.Range("J91").Value = 2023.3334
$ws.Range("K91").Value = 1848
$ws.Range("L91").Value = 2023.3334
$ws.Range("M91").Value = -444
$ws.Range("N91").Value = -4831.3334
$ws.Range("H97").Value = 405.5
$ws.Range("I97").Value = 418.66666
$ws.Range("J97").Value = 287
$ws.Range("K97").Value = 418.66666
$ws.Range("L97").Value = 287
$ws.Range("M97").Value = 77.33334000000002
$ws.Range("N97").Value = -1279
$ws.Range("H102").Value = 1595.1177
$ws.Range("I102").Value = 1416.4073
$ws.Range("J102").Value = 2284.4285
$ws.Range("K102").Value = 1416.4073
$ws.Range("L102").Value = 2284.4285
$ws.Range("M102").Value = 205.5926999999999
$ws.Range("N102").Value = -5528.4285
$ws.Range("H110").Value = 4394.84
$ws.Range("J110").Value = 5027.8
$ws.Range("L110").Value = 5027.8
$ws.Range("N110").Value = -9117.799999999999
$ws.Range("H132").Value = 22733306
$ws.Range("I132").Value = 4164.625
$ws.Range("K132").Value = 12493.875
$ws.Range("M132").Value = -9963.875
$ws.Range("H136").Value = 5999
$ws.Range("I136").Value = 4998.5
$ws.Range("J136").Value = 8000
$ws.Range("K136").Value = 14995.5
$ws.Range("L136").Value = 24000
$ws.Range("M136").Value = -12445.5
$ws.Range("N136").Value = -29100

# Sheet 3: BSM
$ws = $wb.Worksheets.Item(3)
$ws.Range("H86").Value = 2377.2144
$ws.Range("I86").Value = 1988.2
$ws.Range("J86").Value = 3349.75
$ws.Range("K86").Value = 1988.2
$ws.Range("L86").Value = 3349.75
$ws.Range("M86").Value = -865.2
$ws.Range("N86").Value = -5595.75
$ws.Range("H89").Value = 2377.2144
$ws.Range("I89").Value = 1988.2
$ws.Range("J89").Value = 3349.75
$ws.Range("K89").Value = 9941
$ws.Range("L89").Value = 16748.75
$ws.Range("M89").Value = -4325
$ws.Range("N89").Value = -27980.75
$ws.Range("H107").Value = 18170.572
$ws.Range("I107").Value = 20832.666
$ws.Range("J107").Value = 2198
$ws.Range("K107").Value = 20832.666
$ws.Range("L107").Value = 2198
$ws.Range("M107").Value = -18912.666
$ws.Range("N107").Value = -6038
$ws.Range("H134").Value = 1995.3636
$ws.Range("I134").Value = 1795
$ws.Range("J134").Value = 3999
$ws.Range("K134").Value = 5385
$ws.Range("L134").Value = 11997
$ws.Range("M134").Value = -2850
$ws.Range("N134").Value = -17067

# Sheet 4: CRP
$ws = $wb.Worksheets.Item(4)
$ws.Range("H10").Value = 1115932.9
$ws.Range("I10").Value = 1669249.6
$ws.Range("J10").Value = 9299.333000000001
$ws.Range("K10").Value = 1669249.6
$ws.Range("L10").Value = 9299.333000000001
$ws.Range("M10").Value = -1669110.6
$ws.Range("N10").Value = -9577.333000000001
$ws.Range("H16").Value = 1149.125
$ws.Range("I16").Value = 1240.5
$ws.Range("K16").Value = 1240.5
$ws.Range("M16").Value = -953.5
$ws.Range("H22").Value = 1345.4445
$ws.Range("I22").Value = 399.25
$ws.Range("K22").Value = 399.25
$ws.Range("M22").Value = -49.25
$ws.Range("H31").Value = 2420.25
$ws.Range("I31").Value = 3036.182
$ws.Range("J31").Value = 2021.7059
$ws.Range("K31").Value = 3036.182
$ws.Range("L31").Value = 2021.7059
$ws.Range("M31").Value = -2741.182
$ws.Range("N31").Value = -2611.7059
$ws.Range("H34").Value = 2420.25
$ws.Range("I34").Value = 3036.182
$ws.Range("J34").Value = 2021.7059
$ws.Range("K34").Value = 3036.182
$ws.Range("L34").Value = 2021.7059
$ws.Range("M34").Value = -2834.182
$ws.Range("N34").Value = -2425.7059
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()
$ws.Range("H59").Value = 45401.066
$ws.Range("J59").Value = 45401.066
$ws.Range("L59").Value = 45401.066
$ws.Range("N59").Value = -47691.066
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("H86").Value = 4481.5713
$ws.Range("I86").Value = 4329.5
$ws.Range("K86").Value = 4329.5
$ws.Range("M86").Value = -3206.5
$ws.Range("H89").Value = 4481.5713
$ws.Range("I89").Value = 4329.5
$ws.Range("K89").Value = 21647.5
$ws.Range("M89").Value = -16031.5
$ws.Range("H99").Value = 3400
$ws.Range("I99").Value = 2801.8
$ws.Range("J99").Value = 3998.2
$ws.Range("K99").Value = 2801.8
$ws.Range("L99").Value = 3998.2
$ws.Range("M99").Value = -1303.8
$ws.Range("N99").Value = -6994.2
$ws.Range("H113").Value = 1149.125
$ws.Range("I113").Value = 1240.5
$ws.Range("K113").Value = 1240.5
$ws.Range("M113").Value = 929.5
$ws.Range("H126").Value = 3400
$ws.Range("I126").Value = 2801.8
$ws.Range("J126").Value = 3998.2
$ws.Range("K126").Value = 8405.400000000001
$ws.Range("L126").Value = 11994.6
$ws.Range("M126").Value = -5935.400000000001
$ws.Range("N126").Value = -16934.6
$ws.Range("H132").Value = 3464.5
$ws.Range("I132").Value = 3548.4119
$ws.Range("K132").Value = 10645.2357
$ws.Range("M132").Value = -8115.235700000001
$ws.Range("H134").Value = 1804.4445
$ws.Range("I134").Value = 1855
$ws.Range("K134").Value = 5565
$ws.Range("M134").Value = -3030

# Sheet 5: CUL
$ws = $wb.Worksheets.Item(5)
$ws.Range("H4").Value = 325187.1
$ws.Range("I4").Value = 421645.62
$ws.Range("J4").Value = 203006.27
$ws.Range("K4").Value = 1264936.86
$ws.Range("L4").Value = 609018.8099999999
$ws.Range("M4").Value = -1264824.86
$ws.Range("N4").Value = -609242.8099999999
$ws.Range("H37").Value = 118357.8
$ws.Range("J37").Value = 118357.8
$ws.Range("L37").Value = 355073.4
$ws.Range("N37").Value = -355297.4
$ws.Range("H39").Value = 6243.375
$ws.Range("J39").Value = 7559.3335
$ws.Range("L39").Value = 22678.0005
$ws.Range("N39").Value = -23266.0005
$ws.Range("H40").Value = 138
$ws.Range("I40").Value = 165.375
$ws.Range("J40").Value = 101.5
$ws.Range("K40").Value = 661.5
$ws.Range("L40").Value = 406
$ws.Range("M40").Value = -592.5
$ws.Range("N40").Value = -544
$ws.Range("H51").Value = 2728.6667
$ws.Range("I51").Value = 4301.6665
$ws.Range("K51").Value = 12904.9995
$ws.Range("M51").Value = -12444.9995
$ws.Range("H68").Value = 1948.7142
$ws.Range("I68").Value = 744.5
$ws.Range("J68").Value = 2149.4167
$ws.Range("K68").Value = 2233.5
$ws.Range("L68").Value = 6448.250100000001
$ws.Range("M68").Value = -1422.5
$ws.Range("N68").Value = -8070.250100000001
$ws.Range("H71").Value = 1948.7142
$ws.Range("I71").Value = 744.5
$ws.Range("J71").Value = 2149.4167
$ws.Range("K71").Value = 6700.5
$ws.Range("L71").Value = 19344.7503
$ws.Range("M71").Value = -2644.5
$ws.Range("N71").Value = -27456.7503
$ws.Range("H107").Value = 1688.5714
$ws.Range("I107").Value = 220.8
$ws.Range("J107").Value = 2007.6522
$ws.Range("K107").Value = 662.4000000000001
$ws.Range("L107").Value = 6022.9566
$ws.Range("M107").Value = 1257.6
$ws.Range("N107").Value = -9862.9566
$ws.Range("H113").Value = 550.4
$ws.Range("I113").Value = 564.5833
$ws.Range("J113").Value = 537.3077
$ws.Range("K113").Value = 1693.7499
$ws.Range("L113").Value = 1611.9231
$ws.Range("M113").Value = 476.2501
$ws.Range("N113").Value = -5951.9231
$ws.Range("H121").Value = 809656
$ws.Range("I121").Value = 100416.7
$ws.Range("J121").Value = 1252930.5
$ws.Range("K121").Value = 301250.1
$ws.Range("L121").Value = 3758791.5
$ws.Range("M121").Value = -299940.1
$ws.Range("N121").Value = -3761411.5
$ws.Range("H131").Value = 72944.28999999999
$ws.Range("I131").Value = 1833.3334
$ws.Range("J131").Value = 92338.17999999999
$ws.Range("K131").Value = 5500.0002
$ws.Range("L131").Value = 277014.54
$ws.Range("M131").Value = -460.0002000000004
$ws.Range("N131").Value = -287094.54
$ws.Range("H137").Value = 1143726.8
$ws.Range("I137").Value = 126521.75
$ws.Range("J137").Value = 2500000
$ws.Range("K137").Value = 379565.25
$ws.Range("L137").Value = 7500000
$ws.Range("M137").Value = -374465.25
$ws.Range("N137").Value = -7510200

# Sheet 6: GSM
$ws = $wb.Worksheets.Item(6)
$ws.Range("H40").Value = 34999
$ws.Range("J40").Value = 34999
$ws.Range("L40").Value = 34999
$ws.Range("N40").Value = -35301
$ws.Range("H80").Value = 2952
$ws.Range("I80").Value = 2220
$ws.Range("J80").Value = 3074
$ws.Range("K80").Value = 2220
$ws.Range("L80").Value = 3074
$ws.Range("M80").Value = -1222
$ws.Range("N80").Value = -5070
$ws.Range("H83").Value = 2952
$ws.Range("I83").Value = 2220
$ws.Range("J83").Value = 3074
$ws.Range("K83").Value = 11100
$ws.Range("L83").Value = 15370
$ws.Range("M83").Value = -6108
$ws.Range("N83").Value = -25354
$ws.Range("H102").Value = 2412.6667
$ws.Range("I102").Value = 2260.375
$ws.Range("K102").Value = 2260.375
$ws.Range("M102").Value = -638.375
$ws.Range("H107").Value = 739.38464
$ws.Range("J107").Value = 628.44446
$ws.Range("L107").Value = 628.44446
$ws.Range("N107").Value = -4468.44446
$ws.Range("H113").Value = 3135.5
$ws.Range("I113").Value = 3135.5
$ws.Range("K113").Value = 3135.5
$ws.Range("M113").Value = -965.5
$ws.Range("H132").Value = 4334.1665
$ws.Range("I132").Value = 4334.1665
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 13002.4995
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -10472.4995
$ws.Range("N132").ClearContents()

# Sheet 7: LTW
$ws = $wb.Worksheets.Item(7)
$ws.Range("H2").Value = 6843.35
$ws.Range("J2").Value = 15827.833
$ws.Range("L2").Value = 15827.833
$ws.Range("N2").Value = -16051.833
$ws.Range("H22").Value = 7028.75
$ws.Range("I22").Value = 1000
$ws.Range("K22").Value = 1000
$ws.Range("M22").Value = -705
$ws.Range("H27").Value = 7028.75
$ws.Range("I27").Value = 1000
$ws.Range("K27").Value = 1000
$ws.Range("M27").Value = -893
$ws.Range("H40").Value = 3595.56
$ws.Range("I40").Value = 3289.476
$ws.Range("J40").Value = 5202.5
$ws.Range("K40").Value = 3289.476
$ws.Range("L40").Value = 5202.5
$ws.Range("M40").Value = -3153.476
$ws.Range("N40").Value = -5474.5
$ws.Range("H46").Value = 7400.4
$ws.Range("J46").Value = 7500.5
$ws.Range("L46").Value = 7500.5
$ws.Range("N46").Value = -7876.5
$ws.Range("H61").Value = 2565.5715
$ws.Range("I61").Value = 2642.4736
$ws.Range("J61").Value = 1835
$ws.Range("K61").Value = 2642.4736
$ws.Range("L61").Value = 1835
$ws.Range("M61").Value = -2440.4736
$ws.Range("N61").Value = -2239
$ws.Range("H68").Value = 5736.6924
$ws.Range("I68").Value = 4430
$ws.Range("J68").Value = 10092.333
$ws.Range("K68").Value = 4430
$ws.Range("L68").Value = 10092.333
$ws.Range("M68").Value = -3681
$ws.Range("N68").Value = -11590.333
$ws.Range("H71").Value = 5736.6924
$ws.Range("I71").Value = 4430
$ws.Range("J71").Value = 10092.333
$ws.Range("K71").Value = 22150
$ws.Range("L71").Value = 50461.665
$ws.Range("M71").Value = -18406
$ws.Range("N71").Value = -57949.665
$ws.Range("H113").Value = 2565.5715
$ws.Range("I113").Value = 2642.4736
$ws.Range("J113").Value = 1835
$ws.Range("K113").Value = 2642.4736
$ws.Range("L113").Value = 1835
$ws.Range("M113").Value = -472.4735999999998
$ws.Range("N113").Value = -6175
$ws.Range("H122").Value = 7551.7144
$ws.Range("I122").Value = 7266.278
$ws.Range("K122").Value = 21798.834
$ws.Range("M122").Value = -19348.834
$ws.Range("H136").Value = 43481748
$ws.Range("I136").Value = 3158.7368
$ws.Range("J136").Value = 250005040
$ws.Range("K136").Value = 9476.2104
$ws.Range("L136").Value = 750015120
$ws.Range("M136").Value = -6926.2104
$ws.Range("N136").Value = -750020220

# Sheet 8: WVR
$ws = $wb.Worksheets.Item(8)
$ws.Range("H69").Value = 41394.125
$ws.Range("I69").Value = 28886
$ws.Range("K69").Value = 28886
$ws.Range("M69").Value = -28137
$ws.Range("H72").Value = 41394.125
$ws.Range("I72").Value = 28886
$ws.Range("K72").Value = 86658
$ws.Range("M72").Value = -82914
$ws.Range("H107").Value = 2437.2693
$ws.Range("I107").Value = 2042.3158
$ws.Range("J107").Value = 3509.2856
$ws.Range("K107").Value = 6126.9474
$ws.Range("L107").Value = 10527.8568
$ws.Range("M107").Value = -4206.9474
$ws.Range("N107").Value = -14367.8568
$ws.Range("H113").Value = 2166.4666
$ws.Range("I113").Value = 602
$ws.Range("K113").Value = 1806
$ws.Range("M113").Value = 364
$ws.Range("H122").Value = 2187.5293
$ws.Range("I122").Value = 1968.9
$ws.Range("K122").Value = 5906.700000000001
$ws.Range("M122").Value = -3456.700000000001
$ws.Range("H132").Value = 2328.1538
$ws.Range("I132").Value = 2440.5652
$ws.Range("J132").Value = 1466.3334
$ws.Range("K132").Value = 7321.6956
$ws.Range("L132").Value = 4399.0002
$ws.Range("M132").Value = -4791.6956
$ws.Range("N132").Value = -9459.0002
$ws.Range("H136").Value = 2467.04
$ws.Range("I136").Value = 2403.5833
$ws.Range("K136").Value = 7210.749899999999
$ws.Range("M136").Value = -4660.749899999999
